$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the checklist status cells (dropdown list validated: "✅,❌")
$ws.Range("B9").Value = "✅"
$ws.Range("B11").Value = "❌"
$ws.Range("B14").Value = "✅"
$ws.Range("B15").Value = "✅"

# Update the active selection to reflect where the user left off editing
$ws.Range("I16").Select()
